$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'66.153.70"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +2.36%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'3.239.54"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +5.50%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.03%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'574.93"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +2.36%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'152.49"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +6.03%  "
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'3.227.18"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +5.41%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.513"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +3.43%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'7.05"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +8.50%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.163"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +3.30%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.486"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +2.90%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'37.52"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +2.19%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'0.0000233"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +3.38%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'3.762.34"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +5.58%  "
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'555.41"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +10.94%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'66.194.96"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +2.35%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'3.248.28"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  +5.81%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +2.48%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'7.07"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +4.34%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'14.35"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +3.21%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.740"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +6.15%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'7.73"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +6.32%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'13.52"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  +4.56%  "
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = "'81.57"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +2.46%  "
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.20%  "
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'9.24"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +16.42%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'2.93"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  +4.75%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'2.22"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  +4.25%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'27.63"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +4.61%  "
$c.Style = "Normal"
$c = $ws.Range("B31")
$c.Value = "'FirstDigitalUSD"
$c.Style = "Normal"
$c = $ws.Range("C31")
$c.Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  +0.12%  "
$c.Style = "Normal"
$c = $ws.Range("B32")
$c.Value = "'Stacks"
$c.Style = "Normal"
$c = $ws.Range("C32")
$c.Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'2.72"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +1.23%  "
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +4.90%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'559.50"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +8.02%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'5.68"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +2.15%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'6.34"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +4.44%  "
$c.Style = "Normal"
$c = $ws.Range("B37")
$c.Value = "'OKB"
$c.Style = "Normal"
$c = $ws.Range("C37")
$c.Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'55.26"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +2.70%  "
$c.Style = "Normal"
$c = $ws.Range("B38")
$c.Value = "'VeChain"
$c.Style = "Normal"
$c = $ws.Range("C38")
$c.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'0.0452"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +10.11%  "
$c.Style = "Normal"
$c = $ws.Range("B39")
$c.Value = "'Hedera"
$c.Style = "Normal"
$c = $ws.Range("C39")
$c.Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.0857"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  +5.38%  "
$c.Style = "Normal"
$c = $ws.Range("B40")
$c.Value = "'Kaspa"
$c.Style = "Normal"
$c = $ws.Range("C40")
$c.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'0.130"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +4.11%  "
$c.Style = "Normal"
$c = $ws.Range("B41")
$c.Value = "'dogwifhat"
$c.Style = "Normal"
$c = $ws.Range("C41")
$c.Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'3.02"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +10.28%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'3.145.09"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  +6.57%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'8.56"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +0.70%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.273"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +9.05%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'2.27"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +5.20%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'26.33"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  +2.61%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.999"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'0.0₃0549"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -0.12%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'124.46"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  +2.46%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.112"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +1.01%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'2.23"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  +5.91%  "
$c.Style = "Normal"
